$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 195
$ws.Range("I12").Value = 190
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 190
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -20
$ws.Range("N12").Value = -540
$ws.Range("H21").Value = 25500
$ws.Range("I21").Value = 25500
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 25500
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -25032
$ws.Range("H23").Value = 25500
$ws.Range("I23").Value = 25500
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 25500
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -25266
$ws.Range("H100").Value = 2604.5
$ws.Range("I100").Value = 1889.8889
$ws.Range("J100").Value = 3033.2666
$ws.Range("K100").Value = 1889.8889
$ws.Range("L100").Value = 3033.2666
$ws.Range("M100").Value = -1348.8889
$ws.Range("N100").Value = -4115.2666
$ws.Range("H137").Value = 1396.4667
$ws.Range("I137").Value = 1062.8438
$ws.Range("J137").Value = 2217.6924
$ws.Range("K137").Value = 3188.5314
$ws.Range("L137").Value = 6653.0772
$ws.Range("M137").Value = -638.5314000000003
$ws.Range("N137").Value = -11753.0772
$ws.Range("H138").Value = 3760.426
$ws.Range("I138").Value = 1068.5526
$ws.Range("J138").Value = 10153.625
$ws.Range("K138").Value = 3205.6578
$ws.Range("L138").Value = 30460.875
$ws.Range("M138").Value = 1934.3422
$ws.Range("N138").Value = -40740.875
$ws.Range("H139").Value = 53348
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 53348
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 53348
$ws.Range("N139").Value = -63628
$ws.Range("H141").Value = 4451.841
$ws.Range("I141").Value = 2208.9524
$ws.Range("J141").Value = 51552.5
$ws.Range("K141").Value = 6626.8572
$ws.Range("L141").Value = 154657.5
$ws.Range("M141").Value = -1446.8572
$ws.Range("N141").Value = -165017.5
$ws.Range("N21").ClearContents()
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12173.275
$ws.Range("I32").Value = 13448.565
$ws.Range("J32").Value = 7284.6665
$ws.Range("K32").Value = 13448.565
$ws.Range("L32").Value = 7284.6665
$ws.Range("M32").Value = -13161.565
$ws.Range("H45").Value = 1098.75
$ws.Range("I45").Value = 1041.4286
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1041.4286
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -664.4286
$ws.Range("H60").Value = 28683.666
$ws.Range("I60").Value = 28025.5
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 28025.5
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -27292.5
$ws.Range("N60").Value = -31466
$ws.Range("H74").Value = 1126.7354
$ws.Range("I74").Value = 964.6799999999999
$ws.Range("J74").Value = 1576.8889
$ws.Range("K74").Value = 964.6799999999999
$ws.Range("L74").Value = 1576.8889
$ws.Range("M74").Value = -90.67999999999995
$ws.Range("N74").Value = -3324.8889
$ws.Range("H77").Value = 1126.7354
$ws.Range("I77").Value = 964.6799999999999
$ws.Range("J77").Value = 1576.8889
$ws.Range("K77").Value = 4823.4
$ws.Range("L77").Value = 7884.4445
$ws.Range("M77").Value = -455.3999999999996
$ws.Range("N77").Value = -16620.4445
$ws.Range("H97").Value = 909.8823
$ws.Range("I97").Value = 943.8
$ws.Range("J97").Value = 861.4286
$ws.Range("K97").Value = 943.8
$ws.Range("L97").Value = 861.4286
$ws.Range("M97").Value = -447.8
$ws.Range("N97").Value = -1853.4286
$ws.Range("H102").Value = 168803.5
$ws.Range("I102").Value = 2452.5
$ws.Range("J102").Value = 501505.5
$ws.Range("K102").Value = 2452.5
$ws.Range("L102").Value = 501505.5
$ws.Range("M102").Value = -830.5
$ws.Range("N102").Value = -504749.5
$ws.Range("H109").Value = 23361
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 23361
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 23361
$ws.Range("N109").Value = -26135

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2480
$ws.Range("I99").Value = 960
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 960
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = 538
$ws.Range("N99").Value = -6996
$ws.Range("H103").Value = 60000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 60000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 60000
$ws.Range("N103").Value = -62344
$ws.Range("H105").Value = 3313.8
$ws.Range("I105").Value = 3200.5454
$ws.Range("J105").Value = 3625.25
$ws.Range("K105").Value = 3200.5454
$ws.Range("L105").Value = 3625.25
$ws.Range("M105").Value = -1453.5454
$ws.Range("N105").Value = -7119.25
$ws.Range("H134").Value = 1662.7084
$ws.Range("I134").Value = 1506.289
$ws.Range("J134").Value = 4009
$ws.Range("K134").Value = 4518.867
$ws.Range("L134").Value = 12027
$ws.Range("M134").Value = -1983.867
$ws.Range("N134").Value = -17097
$ws.Range("H138").Value = 53362.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 53362.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 53362.5
$ws.Range("N138").Value = -63642.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1163.5358
$ws.Range("I31").Value = 995.2778
$ws.Range("J31").Value = 5706.5
$ws.Range("K31").Value = 995.2778
$ws.Range("L31").Value = 5706.5
$ws.Range("M31").Value = -700.2778
$ws.Range("N31").Value = -6296.5
$ws.Range("H34").Value = 1163.5358
$ws.Range("I34").Value = 995.2778
$ws.Range("J34").Value = 5706.5
$ws.Range("K34").Value = 995.2778
$ws.Range("L34").Value = 5706.5
$ws.Range("M34").Value = -793.2778
$ws.Range("N34").Value = -6110.5
$ws.Range("H132").Value = 323270
$ws.Range("I132").Value = 467184.2
$ws.Range("J132").Value = 2230.6924
$ws.Range("K132").Value = 1401552.6
$ws.Range("L132").Value = 6692.0772
$ws.Range("M132").Value = -1399022.6
$ws.Range("H134").Value = 1455.3508
$ws.Range("I134").Value = 1311.7609
$ws.Range("J134").Value = 2055.818
$ws.Range("K134").Value = 3935.2827
$ws.Range("L134").Value = 6167.454000000001
$ws.Range("M134").Value = -1400.2827
$ws.Range("N134").Value = -11237.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = -731
$ws.Range("H56").Value = 9067.888999999999
$ws.Range("I56").Value = 9067.888999999999
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 9067.888999999999
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8537.888999999999
$ws.Range("H122").Value = 899.8421
$ws.Range("I122").Value = 662.9091
$ws.Range("J122").Value = 1225.625
$ws.Range("K122").Value = 5966.1819
$ws.Range("L122").Value = 11030.625
$ws.Range("M122").Value = -3516.1819
$ws.Range("N122").Value = -15930.625
$ws.Range("H131").Value = 45461636
$ws.Range("I131").Value = 13936
$ws.Range("J131").Value = 71431750
$ws.Range("K131").Value = 41808
$ws.Range("L131").Value = 214295250
$ws.Range("M131").Value = -36768
$ws.Range("N131").Value = -214305330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 13500
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 17666.666
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 17666.666
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -18002.666
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H24").Value = 5798.2
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 6997.75
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 6997.75
$ws.Range("M24").Value = -827
$ws.Range("N24").Value = -7343.75
$ws.Range("H97").Value = 189681.25
$ws.Range("I97").Value = 102770
$ws.Range("J97").Value = 334533.34
$ws.Range("K97").Value = 102770
$ws.Range("L97").Value = 334533.34
$ws.Range("M97").Value = -102274
$ws.Range("H132").Value = 1486.4286
$ws.Range("I132").Value = 1013.4286
$ws.Range("J132").Value = 2905.4285
$ws.Range("K132").Value = 3040.2858
$ws.Range("L132").Value = 8716.2855
$ws.Range("M132").Value = -510.2857999999997
$ws.Range("N132").Value = -13776.2855
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2427.3572
$ws.Range("I16").Value = 2412.0908
$ws.Range("J16").Value = 2483.3333
$ws.Range("K16").Value = 2412.0908
$ws.Range("L16").Value = 2483.3333
$ws.Range("M16").Value = -2242.0908
$ws.Range("N16").Value = -2823.3333
$ws.Range("H20").Value = 8610857
$ws.Range("I20").Value = 20033334
$ws.Range("J20").Value = 43999.75
$ws.Range("K20").Value = 20033334
$ws.Range("L20").Value = 43999.75
$ws.Range("M20").Value = -20033108
$ws.Range("N20").Value = -44451.75
$ws.Range("H50").Value = 48646
$ws.Range("I50").Value = 50000
$ws.Range("J50").Value = 47292
$ws.Range("K50").Value = 50000
$ws.Range("L50").Value = 47292
$ws.Range("M50").Value = -49363
$ws.Range("N50").Value = -48566
$ws.Range("H93").Value = 1100.25
$ws.Range("I93").Value = 800.4286
$ws.Range("J93").Value = 1520
$ws.Range("K93").Value = 800.4286
$ws.Range("L93").Value = 1520
$ws.Range("M93").Value = 447.5714
$ws.Range("N93").Value = -4016
$ws.Range("H132").Value = 2836.5938
$ws.Range("I132").Value = 2337.9038
$ws.Range("J132").Value = 4997.5835
$ws.Range("K132").Value = 7013.7114
$ws.Range("L132").Value = 14992.7505
$ws.Range("M132").Value = -4483.7114
$ws.Range("N132").Value = -20052.7505
$ws.Range("H136").Value = 3002.7163
$ws.Range("I136").Value = 2947.9456
$ws.Range("J136").Value = 3253.75
$ws.Range("K136").Value = 8843.836800000001
$ws.Range("L136").Value = 9761.25
$ws.Range("M136").Value = -6293.836800000001
$ws.Range("N136").Value = -14861.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 14754.75
$ws.Range("I28").Value = 50000
$ws.Range("J28").Value = 3006.3333
$ws.Range("K28").Value = 50000
$ws.Range("L28").Value = 3006.3333
$ws.Range("M28").Value = -49652
$ws.Range("N28").Value = -3702.3333
$ws.Range("H31").Value = 26666.666
$ws.Range("H126").Value = 3939.0386
$ws.Range("I126").Value = 4572.2383
$ws.Range("J126").Value = 1279.6
$ws.Range("K126").Value = 13716.7149
$ws.Range("L126").Value = 3838.8
$ws.Range("M126").Value = -11246.7149
$ws.Range("N126").Value = -8778.799999999999
$ws.Range("H132").Value = 1373.1316
$ws.Range("I132").Value = 986.64
$ws.Range("J132").Value = 2116.3845
$ws.Range("K132").Value = 2959.92
$ws.Range("L132").Value = 6349.1535
$ws.Range("M132").Value = -429.9200000000001
$ws.Range("N132").Value = -11409.1535
